$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume (E) columns for changed rows (2-49),
# and the full Coin/Link/Price/Volume for rows 50-51 where BabyDogeCoin
# drops off the list, Algorand moves up, and USDD is newly added.
# Formula with a leading "'" forces literal text so Excel does not
# reinterpret numeric-looking strings (e.g. "1.00", "27.004.80") as numbers.

# Row 2
$ws.Range("D2").Formula = "'27.004.80"
$ws.Range("E2").Formula = "'  +0.47%  "

# Row 3
$ws.Range("D3").Formula = "'1.561.78"
$ws.Range("E3").Formula = "'  +0.63%  "

# Row 4
$ws.Range("D4").Formula = "'1.00"
$ws.Range("E4").Formula = "'  -0.14%  "

# Row 5
$ws.Range("D5").Formula = "'207.47"
$ws.Range("E5").Formula = "'  +0.28%  "

# Row 6
$ws.Range("E6").Formula = "'  +0.96%  "

# Row 7
$ws.Range("E7").Formula = "'  -0.18%  "

# Row 8
$ws.Range("D8").Formula = "'22.14"
$ws.Range("E8").Formula = "'  +2.06%  "

# Row 9
$ws.Range("D9").Formula = "'0.248"
$ws.Range("E9").Formula = "'  +0.12%  "

# Row 10
$ws.Range("E10").Formula = "'  +1.81%  "

# Row 11
$ws.Range("E11").Formula = "'  +0.11%  "

# Row 12
$ws.Range("D12").Formula = "'1.783.79"
$ws.Range("E12").Formula = "'  +0.51%  "

# Row 13
$ws.Range("D13").Formula = "'1.539.58"
$ws.Range("E13").Formula = "'  -0.96%  "

# Row 14
$ws.Range("E14").Formula = "'  +1.09%  "

# Row 15
$ws.Range("E15").Formula = "'  +1.27%  "

# Row 16
$ws.Range("D16").Formula = "'62.09"
$ws.Range("E16").Formula = "'  +0.60%  "

# Row 17
$ws.Range("D17").Formula = "'27.016.42"
$ws.Range("E17").Formula = "'  +0.50%  "

# Row 18
$ws.Range("D18").Formula = "'0.0₃0707"
$ws.Range("E18").Formula = "'  +2.60%  "

# Row 19
$ws.Range("D19").Formula = "'217.18"
$ws.Range("E19").Formula = "'  +0.28%  "

# Row 20
$ws.Range("E20").Formula = "'  +2.14%  "

# Row 21
$ws.Range("E21").Formula = "'  -0.15%  "

# Row 22
$ws.Range("E22").Formula = "'  +1.50%  "

# Row 23
$ws.Range("E23").Formula = "'  +0.77%  "

# Row 24
$ws.Range("E24").Formula = "'  -2.62%  "

# Row 25
$ws.Range("D25").Formula = "'153.26"
$ws.Range("E25").Formula = "'  -0.15%  "

# Row 26
$ws.Range("E26").Formula = "'  +0.13%  "

# Row 28
$ws.Range("E28").Formula = "'  +1.33%  "

# Row 29
$ws.Range("D29").Formula = "'1.00"
$ws.Range("E29").Formula = "'  -0.15%  "

# Row 30
$ws.Range("E30").Formula = "'  +1.10%  "

# Row 31
$ws.Range("E31").Formula = "'  +1.60%  "

# Row 32
$ws.Range("E32").Formula = "'  +0.57%  "

# Row 33
$ws.Range("E33").Formula = "'  +0.39%  "

# Row 34
$ws.Range("D34").Formula = "'3.10"
$ws.Range("E34").Formula = "'  +3.34%  "

# Row 35
$ws.Range("E35").Formula = "'  +2.94%  "

# Row 36
$ws.Range("D36").Formula = "'1.05"
$ws.Range("E36").Formula = "'  +9.74%  "

# Row 37
$ws.Range("E37").Formula = "'  +1.28%  "

# Row 38
$ws.Range("E38").Formula = "'  +0.71%  "

# Row 39
$ws.Range("D39").Formula = "'0.532"
$ws.Range("E39").Formula = "'  +1.82%  "

# Row 40
$ws.Range("D40").Formula = "'0.809"
$ws.Range("E40").Formula = "'  +0.37%  "

# Row 42
$ws.Range("E42").Formula = "'  +3.11%  "

# Row 43
$ws.Range("D43").Formula = "'5.66"
$ws.Range("E43").Formula = "'  +0.42%  "

# Row 44
$ws.Range("D44").Formula = "'0.998"
$ws.Range("E44").Formula = "'  +0.98%  "

# Row 45
$ws.Range("D45").Formula = "'64.98"
$ws.Range("E45").Formula = "'  +2.08%  "

# Row 46
$ws.Range("E46").Formula = "'  +0.06%  "

# Row 47
$ws.Range("D47").Formula = "'1.702.51"
$ws.Range("E47").Formula = "'  +0.81%  "

# Row 48
$ws.Range("D48").Formula = "'87.63"
$ws.Range("E48").Formula = "'  +1.81%  "

# Row 49
$ws.Range("E49").Formula = "'  +0.68%  "

# Row 50
$ws.Range("B50").Formula = "'Algorand"
$ws.Range("C50").Formula = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Formula = "'0.0956"
$ws.Range("E50").Formula = "'  -0.60%  "

# Row 51
$ws.Range("B51").Formula = "'USDD"
$ws.Range("C51").Formula = "'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").Formula = "'1.00"
$ws.Range("E51").Formula = "'  -0.11%  "
